# This script rotates the data of rows 4, 5 and 6 (columns A, B, E, F, G, H, Q, R, AI)
# in a 3-cycle: row4 <- row6(old), row5 <- row4(old), row6 <- row5(old).
# All other columns/cells in these rows remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original ("before") values for the columns that move.
# Value2 is used instead of Value to get the raw numeric/string data
# without any extra COM variant wrapping.
$orig4 = @{
    A  = $ws.Range("A4").Value2
    B  = $ws.Range("B4").Value2
    E  = $ws.Range("E4").Value2
    F  = $ws.Range("F4").Value2
    G  = $ws.Range("G4").Value2
    H  = $ws.Range("H4").Value2
    Q  = $ws.Range("Q4").Value2
    R  = $ws.Range("R4").Value2
    AI = $ws.Range("AI4").Value2
}

$orig5 = @{
    A  = $ws.Range("A5").Value2
    B  = $ws.Range("B5").Value2
    E  = $ws.Range("E5").Value2
    F  = $ws.Range("F5").Value2
    G  = $ws.Range("G5").Value2
    H  = $ws.Range("H5").Value2
    Q  = $ws.Range("Q5").Value2
    R  = $ws.Range("R5").Value2
    AI = $ws.Range("AI5").Value2
}

$orig6 = @{
    A  = $ws.Range("A6").Value2
    B  = $ws.Range("B6").Value2
    E  = $ws.Range("E6").Value2
    F  = $ws.Range("F6").Value2
    G  = $ws.Range("G6").Value2
    H  = $ws.Range("H6").Value2
    Q  = $ws.Range("Q6").Value2
    R  = $ws.Range("R6").Value2
    AI = $ws.Range("AI6").Value2
}

# Row 4 becomes what row 6 was.
$ws.Range("A4").Value  = $orig6.A
$ws.Range("B4").Value  = $orig6.B
$ws.Range("E4").Value  = $orig6.E
$ws.Range("F4").Value  = $orig6.F
$ws.Range("G4").Value  = $orig6.G
$ws.Range("H4").Value  = $orig6.H
$ws.Range("Q4").Value  = $orig6.Q
$ws.Range("R4").Value  = $orig6.R
$ws.Range("AI4").Value = $orig6.AI

# Row 5 becomes what row 4 was.
$ws.Range("A5").Value  = $orig4.A
$ws.Range("B5").Value  = $orig4.B
$ws.Range("E5").Value  = $orig4.E
$ws.Range("F5").Value  = $orig4.F
$ws.Range("G5").Value  = $orig4.G
$ws.Range("H5").Value  = $orig4.H
$ws.Range("Q5").Value  = $orig4.Q
$ws.Range("R5").Value  = $orig4.R
$ws.Range("AI5").Value = $orig4.AI

# Row 6 becomes what row 5 was.
$ws.Range("A6").Value  = $orig5.A
$ws.Range("B6").Value  = $orig5.B
$ws.Range("E6").Value  = $orig5.E
$ws.Range("F6").Value  = $orig5.F
$ws.Range("G6").Value  = $orig5.G
$ws.Range("H6").Value  = $orig5.H
$ws.Range("Q6").Value  = $orig5.Q
$ws.Range("R6").Value  = $orig5.R
$ws.Range("AI6").Value = $orig5.AI
